# fix: data model and list
#
# Sheet "list 6" holds the Specie -> Breed lookup list. A new breed entry,
# "Bengal Cat" (English) / "Bengal" (French), is added right after the
# "Cat" header row (row 4), pushing the existing "European Shorthair" /
# "Mixed-breed" cat rows and all the dog-breed rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list 6")
$null = $ws.Activate()

# Remember the hyperlinks (cell + URL + display text) before we touch the
# sheet, so we can recreate them shifted down after the row insert.
$linkUrl = "https://creativecommons.org/licenses/by-sa/4.0/"
$oldLinkRows = @(6, 7, 11, 8, 9, 10, 12)

# Duplicate row 4 (Specie=Cat header row) and insert the copy above itself.
# This shifts the original row 4 (and everything below it) down to row 5,
# while the new row 4 keeps the Specie/Cat columns (B-E) populated and
# formatted exactly like the row it was copied from.
$null = $ws.Rows.Item(4).Copy()
$null = $ws.Rows.Item(4).Insert()
$excel.CutCopyMode = 0

# Fill in the new breed entry in the newly inserted row.
$ws.Range("F4").Value() = "Bengal Cat"
$ws.Range("G4").Value() = "Bengal"

# Match the look of the other breed rows in columns F/G (copy format only).
$null = $ws.Range("F5").Copy()
$null = $ws.Range("F4").PasteSpecial(-4122)
$null = $ws.Range("G5").Copy()
$null = $ws.Range("G4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rebuild the hyperlinks on column D, shifted down by one row (everything
# from row 4 onward moved down because of the inserted row above).
$null = $ws.Hyperlinks.Delete()
foreach ($oldRow in $oldLinkRows) {
    $newRow = $oldRow + 1
    $cell = $ws.Cells.Item($newRow, 4)
    $null = $ws.Hyperlinks.Add($cell, $linkUrl, "", "", $linkUrl)
}

# Reflect the active selection recorded in the saved workbook.
$null = $ws.Range("G4").Select()
